# Update the "want to go" count (F) and "lowest ticket price" (G) figures
# that changed between crawl snapshots for the exhibition sheet (sheet 1)
# and the "all types" rollup sheet (sheet 4). The performance / local-life
# sheets are unaffected by this refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: exhibitions ---------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)

$wsExpo.Range("F2").Value = 12819
$wsExpo.Range("G2").Value = 70
$wsExpo.Range("F3").Value = 624
$wsExpo.Range("F5").Value = 32
$wsExpo.Range("F6").Value = 318
$wsExpo.Range("F9").Value = 12815
$wsExpo.Range("F12").Value = 5222
$wsExpo.Range("F16").Value = 29
$wsExpo.Range("F17").Value = 1200
$wsExpo.Range("F18").Value = 35
$wsExpo.Range("F19").Value = 133
$wsExpo.Range("F21").Value = 2850
$wsExpo.Range("F23").Value = 1153
$wsExpo.Range("F24").Value = 3622
$wsExpo.Range("F26").Value = 43

# --- Sheet 4: all types (rollup of every category) ---------------------
$wsAll = $wb.Worksheets.Item(4)

$wsAll.Range("F2").Value = 12819
$wsAll.Range("G2").Value = 70
$wsAll.Range("F3").Value = 624
$wsAll.Range("F5").Value = 32
$wsAll.Range("F6").Value = 318
$wsAll.Range("F10").Value = 12816
$wsAll.Range("F13").Value = 5222
$wsAll.Range("F17").Value = 29
$wsAll.Range("F18").Value = 1200
$wsAll.Range("F19").Value = 35
$wsAll.Range("F20").Value = 133
$wsAll.Range("F22").Value = 2850
$wsAll.Range("F25").Value = 1153
$wsAll.Range("F26").Value = 3622
$wsAll.Range("F28").Value = 43
